$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10, columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
# (columns A,B,C,D,F,L are unchanged)

$data = @{
    2 = @{ E=3; G=46.85851166666667; H=140.575535; I=0.1419057303676978; J=0.1419057303676978; K=3; M=9.682562333333333; N=29.047687; O=0.5358521175370563; P=0.5358521175370563; Q=453.7104600597272; R=4083.394140537545; S=0.07604048610817342; T=0.07604048610817343 }
    3 = @{ E=3; G=46.85851166666667; H=140.575535; I=0.1419057303676978; J=0.1419057303676978; K=3; M=4.196433666666667; N=12.589301; O=0.2322389248810544; P=0.2322389248810544; Q=196.6386359278928; R=1769.747723351035; S=0.03295603425505492; T=0.03295603425505494 }
    4 = @{ E=3; G=46.85851166666667; H=140.575535; I=0.1419057303676978; J=0.1419057303676978; K=3; M=4.190471333333334; N=12.571414; O=0.2319089575818893; P=0.2319089575818892; Q=196.3592498618322; R=1767.23324875649; S=0.03290921000446944; T=0.03290921000446945 }
    5 = @{ E=3; G=283.1772156666667; H=849.531647; I=0.8575703363889615; J=0.8575703363889616; K=3; M=9.682562333333333; N=29.047687; O=0.5358521175370563; P=0.5358521175370563; Q=2741.881042072277; R=24676.92937865049; S=0.4595308806909907; T=0.4595308806909908 }
    6 = @{ E=3; G=283.1772156666667; H=849.531647; I=0.8575703363889615; J=0.8575703363889616; K=3; M=4.196433666666667; N=12.589301; O=0.2322389248810544; P=0.2322389248810544; Q=1188.334401456527; R=10695.00961310875; S=0.1991612129328566; T=0.1991612129328566 }
    7 = @{ E=3; G=283.1772156666667; H=849.531647; I=0.8575703363889615; J=0.8575703363889616; K=3; M=4.190471333333334; N=12.571414; O=0.2319089575818893; P=0.2319089575818892; Q=1186.646004504318; R=10679.81404053886; S=0.1988782427651142; T=0.1988782427651142 }
    8 = @{ E=3; G=0.1730073333333333; H=0.519022; I=0.000523933243340694; J=0.000523933243340694; K=3; M=9.682562333333333; N=29.047687; O=0.5358521175370563; P=0.5358521175370563; Q=1.675154289123778; R=15.076388602114; S=0.0002807507378921687; T=0.0002807507378921687 }
    9 = @{ E=3; G=0.1730073333333333; H=0.519022; I=0.000523933243340694; J=0.000523933243340694; K=3; M=4.196433666666667; N=12.589301; O=0.2322389248810544; P=0.2322389248810544; Q=0.7260137981802222; R=6.534124183622001; S=0.0001216776931428866; T=0.0001216776931428866 }
    10 = @{ E=3; G=0.1730073333333333; H=0.519022; I=0.000523933243340694; J=0.000523933243340694; K=3; M=4.190471333333334; N=12.571414; O=0.2319089575818893; P=0.2319089575818892; Q=0.7249822707897778; R=6.524840437108; S=0.0001215048123056387; T=0.0001215048123056387 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
